# Auto-generated edit script replicating the OOXML diff for cryptos.xlsx
# Updates Price (D) and Volume(1h) (E) columns, and swaps two pairs of rows
# (TRON/WrappedBTC at rows 17-18, and WrappedeETH/Dai at rows 25-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.960.20'
$ws.Range('E2').Value = '  -1.49%  '

$ws.Range('D3').Value = '3.467.15'
$ws.Range('E3').Value = '  -0.45%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = '''600.70'
$ws.Range('E5').Value = '  -0.14%  '

$ws.Range('D6').Value = '''142.41'
$ws.Range('E6').Value = '  -3.77%  '

$ws.Range('D7').Value = '3.465.08'
$ws.Range('E7').Value = '  -0.42%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('E9').Value = '  -1.18%  '

$ws.Range('D10').Value = '''8.17'
$ws.Range('E10').Value = '  +6.81%  '

$ws.Range('D11').Value = '''0.134'
$ws.Range('E11').Value = '  -5.38%  '

$ws.Range('E12').Value = '  -2.85%  '

$ws.Range('D13').Value = '4.060.77'
$ws.Range('E13').Value = '  -0.18%  '

$ws.Range('D14').Value = '''0.0000202'
$ws.Range('E14').Value = '  -4.72%  '

$ws.Range('D15').Value = '''30.27'
$ws.Range('E15').Value = '  -2.91%  '

$ws.Range('D16').Value = '3.472.12'
$ws.Range('E16').Value = '  -0.07%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '66.098.96'
$ws.Range('E17').Value = '  -1.15%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '''0.116'
$ws.Range('E18').Value = '  -0.48%  '

$ws.Range('D19').Value = '''10.41'
$ws.Range('E19').Value = '  +2.25%  '

$ws.Range('E20').Value = '  -3.91%  '

$ws.Range('D21').Value = '''14.67'
$ws.Range('E21').Value = '  -3.74%  '

$ws.Range('D22').Value = '''419.20'
$ws.Range('E22').Value = '  -3.48%  '

$ws.Range('D23').Value = '''0.585'
$ws.Range('E23').Value = '  -3.47%  '

$ws.Range('D24').Value = '''77.42'
$ws.Range('E24').Value = '  -2.15%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '''1.00'
$ws.Range('E25').Value = '  -0.03%  '

$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.617.78'
$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').Value = '''0.0000114'
$ws.Range('E27').Value = '  -4.38%  '

$ws.Range('D28').Value = '''9.37'
$ws.Range('E28').Value = '  -4.49%  '

$ws.Range('D29').Value = '''7.92'
$ws.Range('E29').Value = '  -5.27%  '

$ws.Range('D30').Value = '''2.45'
$ws.Range('E30').Value = '  -1.31%  '

$ws.Range('E31').Value = '  +0.03%  '

$ws.Range('E32').Value = '  -3.49%  '

$ws.Range('D33').Value = '''1.47'
$ws.Range('E33').Value = '  -6.90%  '

$ws.Range('D34').Value = '''25.03'
$ws.Range('E34').Value = '  -1.22%  '

$ws.Range('D35').Value = '3.467.34'
$ws.Range('E35').Value = '  -0.07%  '

$ws.Range('E37').Value = '  -5.55%  '

$ws.Range('D38').Value = '''5.52'
$ws.Range('E38').Value = '  -6.46%  '

$ws.Range('D39').Value = '''7.59'
$ws.Range('E39').Value = '  -3.73%  '

$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.19%  '

$ws.Range('D41').Value = '''169.24'
$ws.Range('E41').Value = '  -2.49%  '

$ws.Range('D42').Value = '''0.0864'
$ws.Range('E42').Value = '  -2.16%  '

$ws.Range('D43').Value = '''0.886'
$ws.Range('E43').Value = '  -1.20%  '

$ws.Range('D44').Value = '''5.09'
$ws.Range('E44').Value = '  -5.78%  '

$ws.Range('E45').Value = '  -8.39%  '

$ws.Range('D46').Value = '''45.65'
$ws.Range('E46').Value = '  -1.71%  '

$ws.Range('D47').Value = '''26.17'
$ws.Range('E47').Value = '  -9.22%  '

$ws.Range('D48').Value = '''1.19'
$ws.Range('E48').Value = '  -3.90%  '

$ws.Range('D49').Value = '''7.11'
$ws.Range('E49').Value = '  -4.67%  '

$ws.Range('D50').Value = '''2.33'
$ws.Range('E50').Value = '  -3.37%  '

$ws.Range('D51').Value = '''0.927'
$ws.Range('E51').Value = '  -4.54%  '

